$wb = $excel.ActiveWorkbook

# A sheet that already carries the "header row" / "index column" style
# (bold font + thin border + centered/top-aligned) used throughout this
# workbook's per-quarter sheets — used purely as a formatting donor via
# Range.Copy so the new sheets end up with byte-identical style indices
# instead of the COM shim minting new (slightly different) style objects.
$donor = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1. The current "总计" sheet (sheetId 6) becomes the new "2022-Q1" sheet:
#    rename it in place and replace its contents with the 2022-Q1 fund
#    holdings table.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Seed an 8x8 block (header + 7 data rows) from the donor sheet so the
# header row (B1:H1) and the index column (A2:A8) inherit the shared
# "style 2" formatting, then overwrite every value explicitly below.
$donor.Range("A1:H8").Copy($q1.Range("A1"))
$q1.Range("A1").ClearContents()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# row, A(idx), B(code), C(name), D(scale,text), E(position,text),
# F(ratio,text), G(value - text unless noted), H(rank)
$q1Rows = @(
    @(2, 0, "009312", "新疆前海联合价值优选混合A",       "10.35", "92.42", "4.65", "0.4813", $false, 9),
    @(3, 1, "004693", "新疆前海联合泳隽灵活配置混合A",   "9.08",  "93.74", "4.71", "0.4277", $false, 6),
    @(4, 2, "009313", "新疆前海联合价值优选混合C",       "1.67",  "92.42", "4.65", "0.0777", $false, 9),
    @(5, 3, "970015", "申万宏源红利成长灵活配置混合",     "1.04",  "70.39", "2.51", "0.0261", $false, 6),
    @(6, 4, "930602", "国信价值智选混合型集合资产管理计划", "0.50", "67.38", "4.17", "0.0208", $false, 8),
    @(7, 5, "013903", "国泰君安信息行业混合",             "0.25",  "84.06", "3.98", "0.0100", $false, 2),
    @(8, 6, "007042", "新疆前海联合泳隽灵活配置混合C",   "0.00",  "93.74", "4.71", "0",      $true,  6)
)

foreach ($row in $q1Rows) {
    $r    = $row[0]
    $idx  = $row[1]
    $code = $row[2]
    $name = $row[3]
    $dVal = $row[4]
    $eVal = $row[5]
    $fVal = $row[6]
    $gVal = $row[7]
    $gIsNumber = $row[8]
    $hVal = $row[9]

    $q1.Cells.Item($r, 1).Value = $idx

    $q1.Cells.Item($r, 2).Value = "'" + $code
    $q1.Cells.Item($r, 2).Style = "Normal"

    $q1.Cells.Item($r, 3).Value = $name

    $q1.Cells.Item($r, 4).Value = "'" + $dVal
    $q1.Cells.Item($r, 4).Style = "Normal"

    $q1.Cells.Item($r, 5).Value = "'" + $eVal
    $q1.Cells.Item($r, 5).Style = "Normal"

    $q1.Cells.Item($r, 6).Value = "'" + $fVal
    $q1.Cells.Item($r, 6).Style = "Normal"

    if ($gIsNumber) {
        $q1.Cells.Item($r, 7).Value = [double]$gVal
    } else {
        $q1.Cells.Item($r, 7).Value = "'" + $gVal
        $q1.Cells.Item($r, 7).Style = "Normal"
    }

    $q1.Cells.Item($r, 8).Value = $hVal
}

# ---------------------------------------------------------------------------
# 2. Add a brand new "总计" sheet at the end of the workbook with the
#    previous totals table plus a new leading 2022-Q1 row.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Seed a 7x4 block (header + 6 data rows) from the donor sheet, trimmed to
# columns A-D, for the same style-index reasons as above.
$donor.Range("A1:D7").Copy($total.Range("A1"))
$total.Range("A1").ClearContents()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(2, 0, "2022-Q1", 7,  1.04),
    @(3, 1, "2021-Q4", 20, 5.88),
    @(4, 2, "2021-Q3", 14, 10.35),
    @(5, 3, "2021-Q2", 13, 1.46),
    @(6, 4, "2021-Q1", 33, 16.63),
    @(7, 5, "2020-Q4", 30, 15.55)
)

foreach ($row in $totalRows) {
    $r = $row[0]
    $total.Cells.Item($r, 1).Value = $row[1]
    $total.Cells.Item($r, 2).Value = $row[2]
    $total.Cells.Item($r, 3).Value = $row[3]
    $total.Cells.Item($r, 4).Value = $row[4]
}

Write-Output "edit complete"
